$d = $word.ActiveDocument

# Locate the "Acknowedgement" Heading1 paragraph (the section heading),
# independent of any hard-coded paragraph index.
$findRange = $d.Content
$findRange.Find.Execute("Acknowedgement", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$headingPara = $findRange.Paragraphs(1)
$headingRange = $headingPara.Range

# Clear any implicit/explicit numbering that might be inherited by the
# Heading1 style for this paragraph; this stamps an explicit
# <w:numPr><w:ilvl w:val="0"/><w:numId w:val="0"/></w:numPr> override
# on the paragraph (i.e. "no numbering"), matching Word's own behaviour
# when numbering is removed from a paragraph.
$headingRange.ListFormat.RemoveNumbers()

# The hidden "_GoBack" bookmark (Word's "last edit location" marker) is
# moved from the end of the following body-text paragraph to the very
# start of this heading paragraph - re-adding a bookmark with the same
# name relocates it instead of duplicating it.
$goBackRange = $headingPara.Range
$goBackRange.Collapse(1)
$d.Bookmarks.Add("_GoBack", $goBackRange)
